# Update the "Runs" worksheet: replace rows 1-4 (Run079-Run082) with new
# run identifiers (Run083-Run086) and new result values.
$wb = $excel.ActiveWorkbook

$wsRuns = $wb.Worksheets.Item("Runs")

$runsData = @(
    @("Run083", "0.432 (0.428)", "0.498 (0.128)", "0.408 (0.492)", "0.602 (0.463)", "0.789 (0.164)", "0.560 (0.498)", "0.752 (0.395)", "0.683 (0.164)", "0.765 (0.424)", "0.571 (0.447)", "0.590 (0.184)", "0.566 (0.496)"),
    @("Run084", "0.308 (0.417)", "0.617 (0.203)", "0.229 (0.421)", "0.532 (0.446)", "0.735 (0.186)", "0.437 (0.498)", "0.625 (0.427)", "0.715 (0.171)", "0.585 (0.493)", "0.463 (0.450)", "0.684 (0.193)", "0.386 (0.487)"),
    @("Run085", "0.638 (0.448)", "0.655 (0.107)", "0.635 (0.482)", "0.577 (0.421)", "0.767 (0.120)", "0.443 (0.499)", "0.686 (0.411)", "0.767 (0.146)", "0.653 (0.477)", "0.651 (0.430)", "0.738 (0.140)", "0.622 (0.485)"),
    @("Run086", "0.488 (0.466)", "0.641 (0.236)", "0.452 (0.498)", "0.504 (0.466)", "0.784 (0.119)", "0.427 (0.496)", "0.684 (0.436)", "0.514 (0.178)", "0.707 (0.455)", "0.563 (0.464)", "0.640 (0.221)", "0.547 (0.498)")
)

for ($r = 0; $r -lt $runsData.Length; $r++) {
    $rowValues = $runsData[$r]
    $rowNum = $r + 1
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $colNum = $c + 2  # column B is index 2
        $wsRuns.Cells.Item($rowNum, $colNum).Value = $rowValues[$c]
    }
}

# Update the "Iterations" worksheet: replace row 1 (Iteration23) with the
# new iteration identifier (Iteration24) and new result values.
$wsIterations = $wb.Worksheets.Item("Iterations")

$iterationValues = @("Iteration24", "0.47 (0.12)", "0.60 (0.06)", "0.43 (0.14)", "0.55 (0.04)", "0.77 (0.02)", "0.47 (0.05)", "0.69 (0.05)", "0.67 (0.09)", "0.68 (0.07)", "0.56 (0.07)", "0.66 (0.05)", "0.53 (0.09)")

for ($c = 0; $c -lt $iterationValues.Length; $c++) {
    $colNum = $c + 2  # column B is index 2
    $wsIterations.Cells.Item(1, $colNum).Value = $iterationValues[$c]
}
